# Implement the scenario PDP - Mainland UK Delivery Option
#
# The TestData sheet's D2 cell held a numeric value (318498). The edit
# replaces it with a text value "486442" entered with a leading
# apostrophe (quote-prefix), i.e. a number stored as text, and moves the
# active selection to E6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# D2: replace the plain number with the text "486442" (quote-prefixed,
# so Excel keeps it as text even though it looks numeric).
$ws.Range("D2").Value = "'486442"

# Move/update the active selection to E6.
$ws.Range("E6").Select()
